$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2026-01 (row 26)
$ws.Range("B26").Value = 6498
$ws.Range("C26").Value = 1014
$ws.Range("D26").Value = 6060087
$ws.Range("E26").Value = 932.6080332409972
$ws.Range("F26").Value = 9.819165117458173
$ws.Range("G26").Value = 7.643312101910826
$ws.Range("H26").Value = 26.20151593750242
